# Cotações atualizadas - 2025-11-27
# Adds a new row (83) with the quotation data for 2025-11-27 (Excel date serial 45988),
# matching the style/number format already used in column A for previous rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A83").Value = 45988
$ws.Range("A83").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B83").Value = "21,7347"
$ws.Range("C83").Value = "16,0918"
$ws.Range("D83").Value = "15,5203"
$ws.Range("E83").Value = "15,5203"
